# [IMP] z0bug_odoo: new test data
# Add a new requirement value for the "due_cost_service_id" test row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = "l10n_it_ricevute_bancarie"

# Move the active selection down past the data, mirroring the author's
# final cursor position after editing the sheet.
$ws.Range("A19").Select()
